$d = $word.ActiveDocument

# --- 1. Extend first paragraph's text and append a new red-colored run ---
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
# Trim paragraph mark off the end so we operate on the run text only.
$r1.End = $r1.End - 1
$r1.Text = "This is a Microsoft word document.  "

# Insert the new run right after the existing text, before the paragraph mark.
$newRange = $r1.Duplicate
$newRange.Collapse(0)  # wdCollapseEnd
$newRange.InsertAfter([string]::Format("(This is a change {0} Version for branch alternate)", [char]0x2013))
$newRange.Font.Color = 192  # wdColor long (0x00BBGGRR) for RGB C00000

# --- 2. Mark the "Normal (Web)" style as semi-hidden (hidden from the
#        recommended style list / gallery until used) ---
$style = $d.Styles("Normal (Web)")
$style.Visibility = $false
